$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.755.88"
Set-TextValue "D3" "1.649.51"
Set-TextValue "E3" "  +0.78%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "215.64"
Set-TextValue "E5" "  +1.16%  "
Set-TextValue "E6" "  +0.65%  "
Set-TextValue "E7" "  -0.07%  "
Set-TextValue "E8" "  -0.57%  "
Set-TextValue "E9" "  +0.86%  "
Set-TextValue "D10" "19.44"
Set-TextValue "E10" "  +1.59%  "
Set-TextValue "E11" "  +0.63%  "
Set-TextValue "D12" "1.878.44"
Set-TextValue "E12" "  +0.76%  "
Set-TextValue "D13" "1.690.63"
Set-TextValue "E13" "  +3.17%  "
Set-TextValue "D14" "4.21"
Set-TextValue "E14" "  +2.83%  "
Set-TextValue "E15" "  +1.16%  "
Set-TextValue "D16" "66.32"
Set-TextValue "E16" "  +4.42%  "
Set-TextValue "D17" "26.820.15"
Set-TextValue "E17" "  +0.74%  "
Set-TextValue "D18" "0.0₃0757"
Set-TextValue "E18" "  +1.49%  "
Set-TextValue "D19" "223.65"
Set-TextValue "E19" "  +1.87%  "
Set-TextValue "D21" "4.40"
Set-TextValue "E21" "  +2.16%  "
Set-TextValue "E22" "  +3.03%  "
Set-TextValue "D23" "9.58"
Set-TextValue "E23" "  +0.75%  "
Set-TextValue "D24" "2.16"
Set-TextValue "E24" "  +12.31%  "
Set-TextValue "D25" "148.21"
Set-TextValue "E25" "  -0.79%  "
Set-TextValue "E26" "  +0.04%  "
Set-TextValue "E27" "  -0.63%  "
Set-TextValue "D28" "7.10"
Set-TextValue "E28" "  +3.07%  "
Set-TextValue "D29" "15.99"
Set-TextValue "E29" "  +3.03%  "
Set-TextValue "E30" "  +0.05%  "
Set-TextValue "E31" "  +0.16%  "
Set-TextValue "D32" "3.46"
Set-TextValue "E32" "  +5.02%  "
Set-TextValue "D33" "3.07"
Set-TextValue "E33" "  +4.90%  "
Set-TextValue "D34" "1.302.34"
Set-TextValue "E34" "  +10.23%  "
Set-TextValue "D35" "1.57"
Set-TextValue "E35" "  +4.06%  "
Set-TextValue "D36" "0.0183"
Set-TextValue "E36" "  +5.81%  "
Set-TextValue "D38" "0.825"
Set-TextValue "E38" "  +1.83%  "
Set-TextValue "D39" "0.525"
Set-TextValue "E39" "  +3.12%  "
Set-TextValue "E41" "  +3.01%  "
Set-TextValue "E42" "  -3.08%  "
Set-TextValue "D43" "5.42"
Set-TextValue "E43" "  +0.37%  "
Set-TextValue "D44" "1.789.77"
Set-TextValue "E44" "  +0.94%  "
Set-TextValue "D45" "93.78"
Set-TextValue "D46" "61.10"
Set-TextValue "E46" "  +11.31%  "
Set-TextValue "D47" "1.62"
Set-TextValue "E47" "  +4.82%  "
Set-TextValue "D48" "0.0518"
Set-TextValue "E48" "  +1.02%  "
Set-TextValue "D49" "7.84"
Set-TextValue "E49" "  +1.71%  "
Set-TextValue "D50" "0.0981"
Set-TextValue "E50" "  +3.45%  "
Set-TextValue "E51" "  -0.87%  "
